$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift the dates in A2:A6 forward by 4 days
$ws.Range("A2").Value = 45359
$ws.Range("A3").Value = 45360
$ws.Range("A4").Value = 45361
$ws.Range("A5").Value = 45362
$ws.Range("A6").Value = 45363

# Explicitly (re)apply the General number format to the Employee Name column
$ws.Range("B2:B6").NumberFormat = "General"

# Update the active selection on the sheet
$ws.Range("C10").Select()
